# Updated the way to load pos and neg training images
# Adds a new "Sheet2" (after the existing Sheet1) containing the
# hyper-parameter scaling calculations, and makes it the active sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 so the tab order becomes
# Sheet1, Sheet2 (matches <sheets> order in the diff) and Sheet2 becomes
# the active tab (activeTab="1").
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws2.Name = "Sheet2"

# ---- Block 1 : max_a / b / rand inputs -----------------------------------
$ws2.Range("C2").Value = "max_a"
$ws2.Range("D2").Value = 20
$ws2.Range("F2").Value = "hyper-parameter"
$ws2.Range("J2").Value = "a is max for y=0"

$ws2.Range("C3").Value = "b"
$ws2.Range("D3").Value = 10
$ws2.Range("F3").Value = "hyper-parameter"
$ws2.Range("J3").Value = "max_a=log_{q}(q+x) "

$ws2.Range("C4").Value = "rand"
$ws2.Range("D4").Value = 0.5
$ws2.Range("F4").Value = "random"
$ws2.Range("J4").Value = "'(q+x)=q^max_a"

$ws2.Range("J5").Value = "'x=q^max_a-q"

# ---- Block 2 : derived q / x / y ------------------------------------------
$ws2.Range("C6").Value = "q"
$ws2.Range("D6").Formula = "=1+D3/(D2-1)"
$ws2.Range("F6").Value = "q=1+b/(max_a-1)"

$ws2.Range("C7").Value = "x"
$ws2.Range("D7").Formula = "=D6^D2-D6"
$ws2.Range("F7").Value = "'x=q^max_a-q"
$ws2.Range("J7").Value = "a=log_{q+0,5*(q^max_a-q)}(q+q^max_a-q)"

$ws2.Range("C8").Value = "y"
$ws2.Range("D8").Formula = "=(1-D4)*D7"
$ws2.Range("F8").Value = "y=(1-rand)*x"
$ws2.Range("J8").Value = "a=log_{0,5*(q^max_a+q)}(q^max_a)"

$ws2.Range("J9").Value = "a=ln(q^max_a)/ln(0,5*(q^max_a+q))"

# ---- Block 3 : a / a_0,5 ----------------------------------------------------
$ws2.Range("C10").Value = "a"
$ws2.Range("D10").Formula = "=LOG(D6+D7,2)/LOG(D6+D8,2)"
$ws2.Range("F10").Value = "a=log_{q+y}(q+x)"
$ws2.Range("J10").Value = "a=max_a*ln(q)/(ln(0,5)+ln(q^max_a+q))"

$ws2.Range("C11").Value = "a_0,5"
$ws2.Range("D11").Formula = "=(1+D2)/2"

$ws2.Range("J12").Value = "a -> a': 1 -> 1, max_a -> max_a, a_0,5=(1+max_a)/2 should be reached by rand=0,5"
$ws2.Range("J12").Font.Color = 255

# ---- Notes ------------------------------------------------------------------
$ws2.Range("C19").Value = "TODO maybe construct a method that considers hyper-parameters of the scanning through the image when testing the model"
$ws2.Range("C19").Font.Color = 255

$ws2.Range("D10").Font.Bold = $true

$ws2.Range("C21").Value = "useLowerValues = 100;"
$ws2.Range("C22").Value = "Math.log(useLowerValues*prob + 1) / Math.log(useLowerValues*maxProb + 1))"

# ---- Worked example ---------------------------------------------------------
$ws2.Range("C25").Value = "rand"
$ws2.Range("D25").Value = 0.5

$ws2.Range("C26").Value = "max_a"
$ws2.Range("D26").Value = 10

$ws2.Range("C27").Value = "q"
$ws2.Range("D27").Value = 100

$ws2.Range("C29").Value = "p"
$ws2.Range("D29").Formula = "=LOG10(D27*D25+1)/LOG10(D27+1)"

$ws2.Range("C30").Value = "a"
$ws2.Range("D30").Formula = "=1+(1-D29)*(D26-1)"

# ---- Cosmetics --------------------------------------------------------------
$ws2.Columns.Item(4).ColumnWidth = 9.1666666

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("K29").Select()
